$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "date" cells (A2 / B3) are stored as plain text (shared strings)
# even though the column carries a date number format - this is exactly
# the kind of malformed date string the new validation code needs to
# detect. Assigning the literal text straight to `.Value` would make
# Excel "smart match" it into a real date serial whenever the text
# happens to parse as a valid m/d/yyyy date (e.g. "10/20/2019"), which
# would change the stored type/style. Routing the text through a
# formula + copy + paste-values keeps it as literal text without
# touching number formatting / styles.

# Use a scratch cell far outside the used range as a staging area.
$scratch = $ws.Range("ZZ1")

# B3: "20.03.2020" -> "40/20/2020"
$scratch.Formula = "=""40/20/2020"""
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)

# A2: "20.03.2019" -> "10/20/2019"
$scratch.Formula = "=""10/20/2019"""
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)

# Clean up the scratch cell and clipboard marquee.
$scratch.Value = ""
$excel.CutCopyMode = 0

# Update the view: drop the old scrolled-right viewport (topLeftCell)
# and move the selection to A3.
$ws.Range("A3").Select()
